$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.686.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -6.85%  '
$ws.Range("D3").Value = '''2.547.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.89%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''299.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.53%  '
$ws.Range("D6").Value = '''92.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.72%  '
$ws.Range("E7").Value = '  -3.53%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '''0.548'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.16%  '
$ws.Range("D10").Value = '''35.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.12%  '
$ws.Range("D11").Value = '''0.0804'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.87%  '
$ws.Range("D12").Value = '''7.67'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.94%  '
$ws.Range("E13").Value = '  +5.74%  '
$ws.Range("D14").Value = '''2.941.97'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.63%  '
$ws.Range("D15").Value = '''2.542.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.22%  '
$ws.Range("D16").Value = '''0.876'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.07%  '
$ws.Range("D17").Value = '''14.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.82%  '
$ws.Range("D18").Value = '''42.699.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.74%  '
$ws.Range("D19").Value = '''12.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").Value = '''0.0₃0984'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = '''6.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").Value = '''71.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.23%  '
$ws.Range("D23").Value = '''256.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.46%  '
$ws.Range("E24").Value = '  -3.99%  '
$ws.Range("D25").Value = '''2.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.92%  '
$ws.Range("D26").Value = '''29.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.32%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  -4.20%  '
$ws.Range("D29").Value = '''37.08'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.35%  '
$ws.Range("E30").Value = '  -4.71%  '
$ws.Range("E31").Value = '  -3.82%  '
$ws.Range("D32").Value = '''152.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = '''2.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.04%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''2.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.46%  '
$ws.Range("D35").Value = '''3.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.35%  '
$ws.Range("D36").Value = '''0.0793'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.33%  '
$ws.Range("E37").Value = '  -5.73%  '
$ws.Range("D38").Value = '''24.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.85%  '
$ws.Range("E39").Value = '  -3.94%  '
$ws.Range("D40").Value = '''17.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.16%  '
$ws.Range("E41").Value = '  -4.71%  '
$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D42").Value = '''3.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.58%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''3.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.13%  '
$ws.Range("D44").Value = '''2.082.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.30%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''1.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.68%  '
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").Value = '''84.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.46%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '''9.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("D49").Value = '''2.796.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.81%  '
$ws.Range("D50").Value = '''104.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.44%  '
$ws.Range("D51").Value = '''1.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.10%  '
